# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.858.56'
$ws.Range("E2").Value = '  +2.91%  '

# Row 3
$ws.Range("D3").Value = '2.339.10'
$ws.Range("E3").Value = '  +2.19%  '

# Row 4
$ws.Range("E4").Value = '  +0.89%  '

# Row 5
$ws.Range("D5").Value = '''313.57'
$ws.Range("E5").Value = '  -0.68%  '

# Row 6
$ws.Range("D6").Value = '''108.76'
$ws.Range("E6").Value = '  +4.02%  '

# Row 7
$ws.Range("D7").Value = '''0.632'
$ws.Range("E7").Value = '  +1.35%  '

# Row 8
$ws.Range("E8").Value = '  +0.23%  '

# Row 9
$ws.Range("D9").Value = '''0.619'
$ws.Range("E9").Value = '  +2.79%  '

# Row 10
$ws.Range("D10").Value = '''41.28'
$ws.Range("E10").Value = '  +4.66%  '

# Row 11
$ws.Range("D11").Value = '''0.0921'
$ws.Range("E11").Value = '  +2.02%  '

# Row 12
$ws.Range("D12").Value = '''8.58'
$ws.Range("E12").Value = '  +2.08%  '

# Row 13
$ws.Range("D13").Value = '''1.01'
$ws.Range("E13").Value = '  +1.37%  '

# Row 14
$ws.Range("E14").Value = '  -0.83%  '

# Row 15
$ws.Range("D15").Value = '''15.50'
$ws.Range("E15").Value = '  +1.93%  '

# Row 16
$ws.Range("D16").Value = '2.697.76'
$ws.Range("E16").Value = '  +2.37%  '

# Row 17
$ws.Range("D17").Value = '2.338.15'
$ws.Range("E17").Value = '  +2.54%  '

# Row 18
$ws.Range("D18").Value = '43.633.61'
$ws.Range("E18").Value = '  +2.04%  '

# Row 19
$ws.Range("D19").Value = '''7.58'
$ws.Range("E19").Value = '  +2.03%  '

# Row 20
$ws.Range("E20").Value = '  +1.77%  '

# Row 21
$ws.Range("D21").Value = '''13.08'
$ws.Range("E21").Value = '  -4.04%  '

# Row 22
$ws.Range("D22").Value = '''74.24'
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("D23").Value = '''3.51'
$ws.Range("E23").Value = '  -1.17%  '

# Row 24
$ws.Range("D24").Value = '''268.73'
$ws.Range("E24").Value = '  +2.58%  '

# Row 25
$ws.Range("E25").Value = '  +4.23%  '

# Row 26
$ws.Range("E26").Value = '  -0.32%  '

# Row 27
$ws.Range("D27").Value = '''7.66'
$ws.Range("E27").Value = '  +7.25%  '

# Row 28
$ws.Range("D28").Value = '''11.16'
$ws.Range("E28").Value = '  +2.63%  '

# Row 29
$ws.Range("D29").Value = '''2.30'
$ws.Range("E29").Value = '  +0.43%  '

# Row 30
$ws.Range("D30").Value = '''39.17'
$ws.Range("E30").Value = '  +4.97%  '

# Row 31
$ws.Range("D31").Value = '''22.67'
$ws.Range("E31").Value = '  +1.81%  '

# Row 32
$ws.Range("D32").Value = '''168.26'
$ws.Range("E32").Value = '  +1.00%  '

# Row 33
$ws.Range("D33").Value = '''0.0886'
$ws.Range("E33").Value = '  +1.58%  '

# Row 34
$ws.Range("D34").Value = '''2.78'
$ws.Range("E34").Value = '  +7.46%  '

# Row 35
$ws.Range("E35").Value = '  +1.91%  '

# Row 36
$ws.Range("D36").Value = '''4.77'
$ws.Range("E36").Value = '  +5.14%  '

# Row 37
$ws.Range("E37").Value = '  -0.79%  '

# Row 38
$ws.Range("D38").Value = '''0.0365'
$ws.Range("E38").Value = '  +4.71%  '

# Row 39
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.89'
$ws.Range("E39").Value = '  +8.56%  '

# Row 40
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '''3.81'
$ws.Range("E40").Value = '  -0.17%  '

# Row 41
$ws.Range("D41").Value = '''1.70'
$ws.Range("E41").Value = '  +8.52%  '

# Row 42
$ws.Range("D42").Value = '''104.09'
$ws.Range("E42").Value = '  +13.17%  '

# Row 43
$ws.Range("E43").Value = '  +3.12%  '

# Row 44
$ws.Range("D44").Value = '''71.76'
$ws.Range("E44").Value = '  +2.88%  '

# Row 45
$ws.Range("D45").Value = '''13.29'
$ws.Range("E45").Value = '  +9.07%  '

# Row 46
$ws.Range("E46").Value = '  +0.02%  '

# Row 47
$ws.Range("D47").Value = '''114.03'
$ws.Range("E47").Value = '  +0.49%  '

# Row 48
$ws.Range("D48").Value = '1.669.82'
$ws.Range("E48").Value = '  -3.19%  '

# Row 49
$ws.Range("D49").Value = '''77.35'
$ws.Range("E49").Value = '  -2.22%  '

# Row 50
$ws.Range("D50").Value = '''8.98'
$ws.Range("E50").Value = '  +2.42%  '

# Row 51
$ws.Range("D51").Value = '''1.56'
$ws.Range("E51").Value = '  +11.00%  '
